$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove column R entirely (Budget-derived leftover column no longer used;
# genres/mpaa/distributor data for rows 3 & 4 get rewritten explicitly below).
$ws.Columns.Item(18).Clear()

$ws.Range("O1").Value = 'Distributor'
$ws.Range("P1").Value = 'MPAA'
$ws.Range("Q1").Value = 'Genres'
$ws.Range("H2").Value = '"Everything Everywhere All at Once" is a visually stunning and emotionally resonant film that explores the concept of multiverses in a unique and audacious way. The movie has been praised for its originality, humor, and technical brilliance, with some viewers calling it one of the greatest movies of all time. The film touches on themes of love, family, and the significance of every choice we make. It has been described as a must-watch for everyone, with its blend of action, humor, and heart making it a truly special cinematic experience.'
$ws.Range("I2").Value = 'Based on the reviews, it seems that public opinion on Everything Everywhere All at Once is overwhelmingly positive. Many reviewers praise the film for its originality, emotional depth, and technical brilliance. Some even go as far as to call it one of their favorite movies of all time. The film is also commended for its diverse representation and strong performances from the cast, particularly Jamie Lee Curtis. Overall, it appears that the majority of viewers highly recommend watching Everything Everywhere All at Once.'
$ws.Range("J2").Value = '- audacious - funny - original - technically dazzling - thematically resonant - genetically engineered - multiverse - greatest movie of all time - required viewing - therapy - overwhelming - love - homage - crying - action blockbusters - kinetic - daring - expertly choreographed - endlessly creative - emotional core - intelligent filmmaking - lowbrow humor - multiverses - hot - competing realities - mommy issues - nihilist lesbian representation - pride month'
$ws.Range("O2").Value = 'A24'
$ws.Range("H3").Value = 'Spider-Man: Into the Spider-Verse is hailed as the best Spider-Man film, with stunning animation that viewers want to inject into their veins. It is praised for its understanding of the character of Spider-Man and its creativity. The film reinvigorates the superhero genre and is considered one of the best animated movies ever made. The inclusion of Spider-Man Noir and other unique characters is celebrated, and the film is seen as a groundbreaking achievement in animation. Despite initial skepticism, viewers are blown away by the film''s quality and are eager for a sequel. Overall, Spider-Man: Into the Spider-Verse is a must-watch for fans of all ages.'
$ws.Range("I3").Value = 'Overall, public opinion on Spider-Man: Into the Spider-Verse seems overwhelmingly positive. Reviewers praise the film for its animation, humor, creativity, and understanding of the character of Spider-Man. Many consider it the best Spider-Man movie ever made and one of the best animated movies ever made. The film is also praised for its ability to reinvigorate the superhero genre and appeal to a diverse modern audience. Some reviewers express excitement for a potential sequel and highlight specific elements of the film, such as the voice acting and references to other media. Overall, the film is seen as a groundbreaking achievement in animation and storytelling.'
$ws.Range("J3").Value = '- best spider-man film - creative - wonderfully animated - understanding of character - nice animation - best animated movies ever made - spider-man noir - achievement in animation - comic books textures rhythms colors - spider-noir movie - greatest Spider-Man movie - favorite film of 2018 - animated Bill Sienkiewicz paintings'
$ws.Range("O3").Value = 'Sony Pictures Releasing'
$ws.Range("P3").Value = 'PG'
$ws.Range("Q3").Value = 'Action Adventure Animation Comedy Family Fantasy Sci-Fi'
$ws.Range("H4").Value = 'Inception is a complex heist film dressed in science fiction conventions, following Dom Cobb as he tries to free himself from his past. It is a thought-provoking, layered story with sumptuous aesthetics and a brilliant cast. Some viewers appreciate the film''s depth and craftsmanship, while others find flaws in the dialogue and relationships portrayed. The ambiguous ending leaves audiences questioning the truth and meaning behind the story.'
$ws.Range("I4").Value = 'Based on the reviews, it seems that public opinion on Inception is generally positive. Many people appreciate the complex plot and the attention to detail in the film, as well as the performances of the cast. Some viewers enjoy the themes of the movie and the way it explores the concept of dreams within dreams. However, there are also some criticisms, such as the dialogue and certain character relationships. Overall, it appears that Inception is seen as a visually stunning and thought-provoking film, despite its flaws.'
$ws.Range("J4").Value = 'complex - heist - science fiction - study - man - past - cerebral - pop-masterpiece - thought-provoking - layered - story-telling - sumptuous - aesthetics - flawless - editing - sound design - effects - musical score - pitch-perfect - cast - confident - directorial hand - brilliant - unrivaled - filmmaking - living in your head rent free - chemistry - masterpiece - gay/lesbian solidarity - fanfic - totem - spinning - dreaming - joyous - persona 5 - sexy - bad dialogue - dreams - inception'
$ws.Range("O4").Value = 'Warner Bros.'
$ws.Range("P4").Value = 'PG-13'
$ws.Range("Q4").Value = 'Action Adventure Sci-Fi Thriller'
$ws.Range("H5").Value = 'Spider-Man: Across the Spider-Verse is described as an overwhelming viewing experience, with viewers unable to tear their eyes away from the screen for its entire duration. The animation, humor, soundtrack, and plot are all praised as near perfect, leaving viewers beaming with joy. Some viewers were left speechless and overstimulated after the movie, with one reviewer even considering quitting making live-action films after seeing it. Overall, the film is hailed as a masterpiece, deserving of six stars and mandatory viewing on a giant screen with a full sound system.'
$ws.Range("I5").Value = 'Overall, public opinion on Spider-Man: Across the Spider-Verse seems to be overwhelmingly positive. Reviewers have praised the animation, humor, soundtrack, and overall viewing experience of the film. Many have described it as a mind-blowing and visually dazzling masterpiece, with some even calling it better than the original Spider-Verse movie. Some reviewers have mentioned minor criticisms, such as the movie feeling overstuffed with too many stories, but these do not seem to detract from the overall enjoyment of the film. Overall, the average score of 4.75 reflects the high praise and excitement surrounding Spider-Man: Across the Spider-Verse.'
$ws.Range("J5").Value = 'overwhelming - joyous - colorful - perfect - magical - mind-blowing - incredible - stunning - ambitious - visually dazzling - funny - dazzling - ambitious - heartwarming - unforgettable - psychopath - punk'
$ws.Range("O5").Value = 'Columbia Pictures'
$ws.Range("I6").Value = 'Based on the reviews, it seems that public opinion on Spider-Man: No Way Home is generally positive. Many people were thrilled by the return of Willem Dafoe as Green Goblin and the appearances of Tobey Maguire and Andrew Garfield as other versions of Spider-Man. The emotional moments in the film resonated with audiences, and there were high levels of excitement and engagement in theaters during key scenes. However, some viewers found the film to be overly reliant on CGI and struggled with pacing and exposition in the beginning. Overall, it appears that the film successfully delivered on fan expectations and provided a satisfying experience for many viewers.'
$ws.Range("O6").Value = 'Sony Pictures Releasing'
